# Updates cryptos list values per source refresh (prices/volume%, and three
# coins reshuffled in the ranking: Toncoin/ICP swap rows 32-33; RocketPoolETH/
# Aave/BabyDogeCoin rotate rows 47-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value even when the string looks numeric
    # (e.g. "1.001"), without leaving a lingering custom number format on the
    # cell: remember the existing style, flip to text ("@") for the write, then
    # restore the original style object.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.145.05"
$ws.Range("E2").Value = "  -0.14%  "
Set-TextValue $ws.Range("D3") "1.835.68"
$ws.Range("E3").Value = "  -0.40%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue $ws.Range("D5") "242.31"
$ws.Range("E5").Value = "  +0.63%  "
Set-TextValue $ws.Range("D6") "0.6603"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.11%  "
Set-TextValue $ws.Range("D8") "44.64"
$ws.Range("E8").Value = "  +6.19%  "
Set-TextValue $ws.Range("D9") "0.07407"
$ws.Range("E9").Value = "  -0.17%  "
Set-TextValue $ws.Range("D10") "0.2947"
$ws.Range("E10").Value = "  -0.21%  "
Set-TextValue $ws.Range("D11") "23.13"
$ws.Range("E11").Value = "  +1.34%  "
Set-TextValue $ws.Range("D12") "0.07724"
$ws.Range("E12").Value = "  +0.11%  "
Set-TextValue $ws.Range("D13") "1.844.62"
$ws.Range("E13").Value = "  +0.44%  "
Set-TextValue $ws.Range("D14") "5.001"
$ws.Range("E14").Value = "  -0.04%  "
Set-TextValue $ws.Range("D15") "0.6709"
$ws.Range("E15").Value = "  -0.79%  "
Set-TextValue $ws.Range("D16") "82.51"
$ws.Range("E16").Value = "  -4.09%  "
Set-TextValue $ws.Range("D17") "6.150"
$ws.Range("E17").Value = "  +0.31%  "
Set-TextValue $ws.Range("D18") "0.000008699"
$ws.Range("E18").Value = "  +4.77%  "
Set-TextValue $ws.Range("D19") "29.154.87"
$ws.Range("E19").Value = "  +0.14%  "
Set-TextValue $ws.Range("D20") "2.092.12"
$ws.Range("E20").Value = "  +1.23%  "
Set-TextValue $ws.Range("D21") "12.49"
$ws.Range("E21").Value = "  -0.31%  "
Set-TextValue $ws.Range("D22") "225.68"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue $ws.Range("D24") "7.156"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +0.09%  "
Set-TextValue $ws.Range("D26") "158.53"
$ws.Range("E26").Value = "  -1.23%  "
Set-TextValue $ws.Range("D27") "8.593"
$ws.Range("E27").Value = "  -1.23%  "
Set-TextValue $ws.Range("D28") "0.1389"
Set-TextValue $ws.Range("D29") "18.02"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.24%  "
Set-TextValue $ws.Range("D31") "4.128"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D32") "4.033"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D33") "1.206"
$ws.Range("E33").Value = "  +1.42%  "
Set-TextValue $ws.Range("D34") "0.05395"
$ws.Range("E34").Value = "  +1.40%  "
Set-TextValue $ws.Range("D35") "1.851"
$ws.Range("E35").Value = "  -1.57%  "
Set-TextValue $ws.Range("D36") "0.7457"
$ws.Range("E36").Value = "  -1.91%  "
Set-TextValue $ws.Range("D37") "1.158"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -1.23%  "
Set-TextValue $ws.Range("D39") "1.298.90"
$ws.Range("E39").Value = "  -2.49%  "
Set-TextValue $ws.Range("D40") "0.01794"
$ws.Range("E40").Value = "  -0.29%  "
Set-TextValue $ws.Range("D41") "2.762"
$ws.Range("E41").Value = "  +0.96%  "
Set-TextValue $ws.Range("D42") "6.364"
$ws.Range("E42").Value = "  +6.64%  "
Set-TextValue $ws.Range("D43") "0.9051"
$ws.Range("E43").Value = "  -2.15%  "
Set-TextValue $ws.Range("D44") "1.001"
$ws.Range("E44").Value = "  -0.16%  "
Set-TextValue $ws.Range("D45") "103.52"
$ws.Range("E45").Value = "  +0.09%  "
Set-TextValue $ws.Range("D46") "0.07998"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D47") "0.00000000125"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D48") "1.990.68"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D49") "64.94"
$ws.Range("E49").Value = "  +1.55%  "
Set-TextValue $ws.Range("D50") "0.5136"
$ws.Range("E50").Value = "  -0.53%  "
Set-TextValue $ws.Range("D51") "1.749"
$ws.Range("E51").Value = "  -1.22%  "
